$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1 (VOUT1 block): change set-point resistor divider input B1 from 20 to 12 ---
$ws.Range("B1").Value = 12

# --- New row 8: ADC counts derived from B4 / F4 ---
$ws.Range("A8").Value = "ADC"
$ws.Range("B8").Formula = "=B4/4.096*1024"
$ws.Range("E8").Value = "ADC"
$ws.Range("F8").Formula = "=(F4+2.048)/4.096*1024"

# --- New rows 28-29: Band Gap Voltage / ADC counts ---
$ws.Range("A28").Value = "Band Gap Voltage"
$ws.Range("B28").Value = 1.23
$ws.Range("C28").Value = "V"
$ws.Range("A28").Font.Bold = $true

$ws.Range("A29").Value = "ADC"
$ws.Range("B29").Formula = "=B28/4.096*1024"

# --- Bold the section header labels (they keep their existing fill colors) ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true
$ws.Range("A13").Font.Bold = $true
$ws.Range("E13").Font.Bold = $true
$ws.Range("A22").Font.Bold = $true
$ws.Range("A25").Font.Bold = $true

# --- Column A width widened to fit the longer "Band Gap Voltage" label ---
$ws.Columns("A").ColumnWidth = 19

# --- View state: zoom to 85% and select B23 ---
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("B23").Select() | Out-Null
